# ftlua::push / ftlua::Keys overload table update
# Rewrites the "callGFun / callGTabTFun / callTFun / callTMethod" example
# rows (7-13) to reflect the new ftlua::Keys<ARGS...> based key types, and
# tweaks window/view state (selection, scroll position, column width).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: drop the old standalone "callGFun" example row (A7/C7/G7/H7). ---
# K7 ("callGFun:") is left untouched.
$ws.Range("A7").Clear()
$ws.Range("C7").Clear()
$ws.Range("G7").Clear()
$ws.Range("H7").Clear()

# --- Row 8: was the "callGTabTFun" example, becomes the "callTFun" example. ---
# B8 switches from "T const &tabGlobalKey" (key-cell style) to a plain
# "table" value using the same fill style "table" cells use elsewhere
# (copy it wholesale from B10, which already holds "table" with that style).
$ws.Range("B10").Copy($ws.Range("B8"))
$ws.Range("C8").Value = "ftlua::Keys<ARGS...> const &funTabKeys"
$ws.Range("H8").Value = "callTFun"

# --- Row 9: becomes the relocated "callGFun" example (previously row 7),
# now keyed by ftlua::Keys<ARGS...> instead of a single T const&. ---
$ws.Range("A8").Copy($ws.Range("A9"))
$ws.Range("G8").Copy($ws.Range("G9"))
$ws.Range("G8").Copy($ws.Range("C9"))
$ws.Range("C9").Value = "ftlua::Keys<ARGS...> const &funTabKeys"
$ws.Range("H9").Value = "callGFun"

# --- Row 10: drop the old "callTFun" example row (A10/B10/C10/G10/H10). ---
# J10/K10 ("Lua repr:" / "_G[funGKey](args);") are left untouched.
$ws.Range("A10").Clear()
$ws.Range("B10").Clear()
$ws.Range("C10").Clear()
$ws.Range("G10").Clear()
$ws.Range("H10").Clear()

# --- Row 11: drop the old "callTMethod" example row entirely. ---
$ws.Range("A11").Clear()
$ws.Range("B11").Clear()
$ws.Range("C11").Clear()
$ws.Range("F11").Clear()
$ws.Range("G11").Clear()
$ws.Range("H11").Clear()

# --- Row 12: becomes the relocated "callTMethod" example (previously row 11),
# now keyed by ftlua::Keys<ARGS...> instead of a single T const&. ---
$ws.Range("A13").Copy($ws.Range("A12"))
$ws.Range("B8").Copy($ws.Range("B12"))
$ws.Range("G13").Copy($ws.Range("G12"))
$ws.Range("F13").Copy($ws.Range("F12"))
$ws.Range("G13").Copy($ws.Range("C12"))
$ws.Range("C12").Value = "ftlua::Keys<ARGS...> const &methodTabKeys"
$ws.Range("H12").Value = "callTMethod"

# --- Row 13: "callGTabTMethod" example keeps its shape, only the two key
# cells switch to the ftlua::Keys<ARGS...> spellings. ---
$ws.Range("B13").Value = "ftlua::Keys<ARGS...> const &tabGlobalKeys"
$ws.Range("C13").Value = "ftlua::Keys<ARGS...> const &methodTabKeys"

# --- Column width + view state tweaks ---
$ws.Columns.Item(2).ColumnWidth = 34.33203125
$ws.Range("A3").Select()
$ws.Range("E15").Select()

$wb.Windows.Item(1).Top = 460
